$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a last-modified date serial for each row
# (rows 2-66). The whole column was bulk-updated from 45192 (2023-09-23)
# to 45202 (2023-10-03).
$ws.Range("C2:C66").Value = 45202
